$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 to hold the "Docentes responsaveis" value,
# mirroring the existing label/value row pattern used in row 1/2.
$ws.Rows.Item(13).Insert()

# The freshly inserted row copies formatting from the row above (style on A13);
# clear it completely since the target row has no A cell/style at all.
$ws.Range("A13").Clear()

# Borrow the B/C column formatting (styles 2/3) from row 9 so the new cells
# match the sheet's established per-column styling instead of inheriting
# the row-above style picked up by Insert().
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Objetivos: row content (was incorrectly duplicating the teacher name).
$ws.Range("B10").Value = "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."
$ws.Range("C10").Value = "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."

# newly inserted row holding the Docentes responsaveis value.
$ws.Range("B13").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C13").Value = "5464150 - Mariana Consiglio Kasemodel"

# Programa resumido: row content.
$ws.Range("B14").Value = "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."
$ws.Range("C14").Value = "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."

# Programa: row content.
$ws.Range("B16").Value = "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."
$ws.Range("C16").Value = "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."

# Metodo: row content (was incorrectly duplicating the teacher name).
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Criterio: row content.
$ws.Range("B20").Value = "Média ponderada de atividades e provas."
$ws.Range("C20").Value = "Média ponderada de atividades e provas."

# Norma de recuperacao: row content.
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"

# Bibliografia: row content.
$ws.Range("B22").Value = "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
$ws.Range("C22").Value = "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
